# Actualizacion Datos Personales 4 nov
# Insert a new "1AM" group row (row 2) into the three statistics sheets,
# pushing the existing 1BM..1FM rows down one row.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Estadisticos 1P", "Estadisticos 2P", "Estadisticos Final")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a fresh blank row at row 2, shifting current rows 2-6 down to 3-7.
    $ws.Rows("2:2").Insert()

    # The inserted row inherits formatting from the row above (the bold header
    # with borders); the source data rows use the default (unstyled) format,
    # so strip formatting back to default before filling in values.
    $ws.Range("A2:H2").ClearFormats()

    $ws.Cells.Item(2, 1).Value = "QUÍMICA I"
    $ws.Cells.Item(2, 2).Value = "1AM"

    if ($sheetName -eq "Estadisticos 2P") {
        $ws.Cells.Item(2, 3).Value = 34
        $ws.Cells.Item(2, 4).Value = 34
        $ws.Cells.Item(2, 5).Value = 28
        $ws.Cells.Item(2, 6).Value = 0
        $ws.Cells.Item(2, 7).Value = 0
        # This sheet never uses column H; drop the blank cell Insert() added
        # there so the row matches the source (only A:G populated).
        $ws.Cells.Item(2, 8).Clear()
    } else {
        $ws.Cells.Item(2, 3).Value = 34
        $ws.Cells.Item(2, 4).Value = 6
        $ws.Cells.Item(2, 5).Value = 0
        $ws.Cells.Item(2, 6).Value = 28
        $ws.Cells.Item(2, 7).Value = 82.34999999999999
        $ws.Cells.Item(2, 8).Value = 8.4
    }
}

Write-Host "Inserted 1AM rows into Estadisticos 1P, 2P and Final"
